$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92
$ws.Range("H92").Value = 1096.3846
$ws.Range("I92").Value = 251.42857
$ws.Range("J92").Value = 2082.1667
$ws.Range("K92").Value = 251.42857
$ws.Range("L92").Value = 2082.1667
$ws.Range("M92").Value = 996.57143
$ws.Range("N92").Value = -4578.1667

# Row 98
$ws.Range("H98").Value = 1401
$ws.Range("I98").Value = 1243.4849
$ws.Range("J98").Value = 4000
$ws.Range("K98").Value = 1243.4849
$ws.Range("L98").Value = 4000
$ws.Range("M98").Value = 254.5151000000001
$ws.Range("N98").Value = -6996

# Row 122
$ws.Range("H122").Value = 1401
$ws.Range("I122").Value = 1243.4849
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 3730.4547
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -1280.4547
$ws.Range("N122").Value = -16900

# Row 132
$ws.Range("H132").Value = 4610099.5
$ws.Range("I132").Value = 4610099.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 13830298.5
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -13827768.5

# Row 138
$ws.Range("H138").Value = 1625.027
$ws.Range("I138").Value = 1042.75
$ws.Range("J138").Value = 2700
$ws.Range("K138").Value = 3128.25
$ws.Range("L138").Value = 8100
$ws.Range("M138").Value = 2011.75
$ws.Range("N138").Value = -18380


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 39037.38
$ws.Range("I32").Value = 45766
$ws.Range("J32").Value = 6740
$ws.Range("K32").Value = 45766
$ws.Range("L32").Value = 6740
$ws.Range("M32").Value = -45479
$ws.Range("N32").Value = -7314

# Row 45
$ws.Range("H45").Value = 1276.9
$ws.Range("I45").Value = 856.6667
$ws.Range("J45").Value = 1457
$ws.Range("K45").Value = 856.6667
$ws.Range("L45").Value = 1457
$ws.Range("M45").Value = -479.6667
$ws.Range("N45").Value = -2211

# Row 61
$ws.Range("H61").Value = 2172.647
$ws.Range("I61").Value = 2041.5454
$ws.Range("J61").Value = 2413
$ws.Range("K61").Value = 2041.5454
$ws.Range("L61").Value = 2413
$ws.Range("M61").Value = -1829.5454
$ws.Range("N61").Value = -2837

# Row 74
$ws.Range("H74").Value = 403.94116
$ws.Range("I74").Value = 347.25
$ws.Range("J74").Value = 540
$ws.Range("K74").Value = 347.25
$ws.Range("L74").Value = 540
$ws.Range("M74").Value = 526.75
$ws.Range("N74").Value = -2288

# Row 77
$ws.Range("H77").Value = 403.94116
$ws.Range("I77").Value = 347.25
$ws.Range("J77").Value = 540
$ws.Range("K77").Value = 1736.25
$ws.Range("L77").Value = 2700
$ws.Range("M77").Value = 2631.75
$ws.Range("N77").Value = -11436

# Row 122
$ws.Range("H122").Value = 1924.6471
$ws.Range("I122").Value = 2001.3846
$ws.Range("J122").Value = 1675.25
$ws.Range("K122").Value = 6004.1538
$ws.Range("L122").Value = 5025.75
$ws.Range("M122").Value = -3554.1538
$ws.Range("N122").Value = -9925.75

# Row 132
$ws.Range("H132").Value = 4849.433
$ws.Range("I132").Value = 5899.35
$ws.Range("J132").Value = 2749.6
$ws.Range("K132").Value = 17698.05
$ws.Range("L132").Value = 8248.799999999999
$ws.Range("M132").Value = -15168.05
$ws.Range("N132").Value = -13308.8

# Row 136
$ws.Range("H136").Value = 2172.647
$ws.Range("I136").Value = 2041.5454
$ws.Range("J136").Value = 2413
$ws.Range("K136").Value = 6124.6362
$ws.Range("L136").Value = 7239
$ws.Range("M136").Value = -3574.6362
$ws.Range("N136").Value = -12339


$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = 0

# Row 134
$ws.Range("H134").Value = 52876.3
$ws.Range("I134").Value = 65345.375
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 196036.125
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -193501.125
$ws.Range("N134").Value = -14070


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1596248.9
$ws.Range("I31").Value = 1617.2894
$ws.Range("J31").Value = 4785512
$ws.Range("K31").Value = 1617.2894
$ws.Range("L31").Value = 4785512
$ws.Range("M31").Value = -1322.2894
$ws.Range("N31").Value = -4786102

# Row 34
$ws.Range("H34").Value = 1596248.9
$ws.Range("I34").Value = 1617.2894
$ws.Range("J34").Value = 4785512
$ws.Range("K34").Value = 1617.2894
$ws.Range("L34").Value = 4785512
$ws.Range("M34").Value = -1415.2894
$ws.Range("N34").Value = -4785916

# Row 58
$ws.Range("H58").Value = 1053.4783
$ws.Range("I58").Value = 999.73334
$ws.Range("J58").Value = 1154.25
$ws.Range("K58").Value = 999.73334
$ws.Range("L58").Value = 1154.25
$ws.Range("M58").Value = -796.73334
$ws.Range("N58").Value = -1560.25

# Row 136
$ws.Range("H136").Value = 1053.4783
$ws.Range("I136").Value = 999.73334
$ws.Range("J136").Value = 1154.25
$ws.Range("K136").Value = 2999.20002
$ws.Range("L136").Value = 3462.75
$ws.Range("M136").Value = -449.2000200000002
$ws.Range("N136").Value = -8562.75


$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 2201
$ws.Range("I49").Value = 1401.5
$ws.Range("J49").Value = 3800
$ws.Range("K49").Value = 4204.5
$ws.Range("L49").Value = 11400
$ws.Range("M49").Value = -4048.5
$ws.Range("N49").Value = -11712

# Row 58
$ws.Range("H58").Value = 1205
$ws.Range("I58").Value = 1205
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 3615
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -3487

# Row 131
$ws.Range("H131").Value = 3111.0232
$ws.Range("I131").Value = 14682.857
$ws.Range("J131").Value = 860.94446
$ws.Range("K131").Value = 44048.571
$ws.Range("L131").Value = 2582.83338
$ws.Range("M131").Value = -39008.571
$ws.Range("N131").Value = -12662.83338


$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 8392.333000000001
$ws.Range("I126").Value = 1379.6
$ws.Range("J126").Value = 11089.538
$ws.Range("K126").Value = 4138.799999999999
$ws.Range("L126").Value = 33268.614
$ws.Range("M126").Value = -1668.799999999999
$ws.Range("N126").Value = -38208.614


$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 2372.4644
$ws.Range("I132").Value = 2023
$ws.Range("J132").Value = 3420.8572
$ws.Range("K132").Value = 6069
$ws.Range("L132").Value = 10262.5716
$ws.Range("M132").Value = -3539
$ws.Range("N132").Value = -15322.5716


$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 673.0789
$ws.Range("I132").Value = 562.1667
$ws.Range("J132").Value = 863.2143
$ws.Range("K132").Value = 1686.5001
$ws.Range("L132").Value = 2589.6429
$ws.Range("M132").Value = 843.4999
$ws.Range("N132").Value = -7649.6429

# Row 136
$ws.Range("H136").Value = 1681.25
$ws.Range("I136").Value = 1905
$ws.Range("J136").Value = 786.25
$ws.Range("K136").Value = 5715
$ws.Range("L136").Value = 2358.75
$ws.Range("M136").Value = -3165
$ws.Range("N136").Value = -7458.75

